# Update the dSF column (F) values to reflect repulled data / recalculated mean.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 2
$ws.Range("F7").Value = -5
$ws.Range("F8").Value = 3
$ws.Range("F9").Value = -1
$ws.Range("F10").Value = -6
$ws.Range("F12").Value = -4
